$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so values like "1.00" or "90.652.83" are not
# reinterpreted as numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '90.652.83'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '3.104.94'
$ws.Range("E3").Value = '  -2.29%  '
$ws.Range("D4").Value = '0.996'
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = '232.59'
$ws.Range("E5").Value = '  +6.33%  '
$ws.Range("D6").Value = '627.53'
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").Value = '1.11'
$ws.Range("E7").Value = '  -0.51%  '
$ws.Range("D8").Value = '0.364'
$ws.Range("E8").Value = '  -2.02%  '
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").Value = '3.102.90'
$ws.Range("E10").Value = '  -2.21%  '
$ws.Range("D11").Value = '0.724'
$ws.Range("E11").Value = '  -4.32%  '
$ws.Range("E12").Value = '  -1.57%  '
$ws.Range("D13").Value = '36.61'
$ws.Range("E13").Value = '  +3.90%  '
$ws.Range("D14").Value = '0.0000247'
$ws.Range("E14").Value = '  -0.84%  '
$ws.Range("D15").Value = '5.50'
$ws.Range("E15").Value = '  -2.85%  '
$ws.Range("D16").Value = '90.348.57'
$ws.Range("D18").Value = '3.155.37'
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").Value = '3.87'
$ws.Range("E19").Value = '  +3.10%  '
$ws.Range("B20").Value = 'PEPE'
$ws.Range("C20").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D20").Value = '0.0000212'
$ws.Range("E20").Value = '  -2.84%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = '14.13'
$ws.Range("E21").Value = '  -1.78%  '
$ws.Range("D22").Value = '441.37'
$ws.Range("E22").Value = '  -1.22%  '
$ws.Range("E23").Value = '  +6.51%  '
$ws.Range("D24").Value = '8.95'
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("D25").Value = '5.93'
$ws.Range("E25").Value = '  -1.28%  '
$ws.Range("B26").Value = 'LEO'
$ws.Range("C26").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D26").Value = '7.58'
$ws.Range("E26").Value = '  -1.70%  '
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = '89.34'
$ws.Range("E27").Value = '  +1.67%  '
$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").Value = '12.41'
$ws.Range("E28").Value = '  +0.59%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '3.315.59'
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '9.51'
$ws.Range("E31").Value = '  +2.50%  '
$ws.Range("B32").Value = 'Cronos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D32").Value = '0.159'
$ws.Range("E32").Value = '  -1.97%  '
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D33").Value = '0.980'
$ws.Range("E33").Value = '  -3.47%  '
$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").Value = '0.203'
$ws.Range("E34").Value = '  +16.95%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '26.46'
$ws.Range("E35").Value = '  +1.66%  '
$ws.Range("D36").Value = '3.86'
$ws.Range("E36").Value = '  +3.50%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '0.153'
$ws.Range("E37").Value = '  +5.50%  '
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").Value = '511.80'
$ws.Range("E38").Value = '  -2.94%  '
$ws.Range("B39").Value = 'PancakeSwap'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D39").Value = '1.93'
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").Value = '7.06'
$ws.Range("E40").Value = '  +0.45%  '
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D41").Value = '1.29'
$ws.Range("E41").Value = '  -1.63%  '
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").Value = '0.0904'
$ws.Range("E42").Value = '  +5.79%  '
$ws.Range("B43").Value = 'MantraDAO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D43").Value = '3.57'
$ws.Range("E43").Value = '  +56.50%  '
$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").Value = '0.412'
$ws.Range("E44").Value = '  -0.77%  '
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").Value = '22.19'
$ws.Range("E45").Value = '  -0.10%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = '1.92'
$ws.Range("E47").Value = '  -1.97%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '150.54'
$ws.Range("E48").Value = '  +1.41%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").Value = '0.694'
$ws.Range("E49").Value = '  +5.68%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").Value = '45.11'
$ws.Range("E50").Value = '  +2.33%  '
$ws.Range("B51").Value = 'ImmutableX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D51").Value = '1.35'
$ws.Range("E51").Value = '  +0.19%  '
